# edit.ps1 -- apply the "24.004.42 / Fri Feb 24 03:29:40 UTC 2023" cryptos refresh
# Updates Price (col D) and Volume(1h) (col E) for each coin row, and fixes the
# Algorand / FraxShare row order (rows 38-39 swapped, values refreshed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.004.42"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.651.80"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'310.06"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.3901"
$ws.Range("D8").Value = "'0.3825"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "'1.350"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'0.08451"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "'23.85"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "'7.082"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "'8.012"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "'0.00001311"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").Value = "1.649.60"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'94.61"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "'0.06988"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'19.67"
$ws.Range("E20").Value = "  -4.50%  "
$ws.Range("D21").Value = "'6.983"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'13.79"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "23.995.19"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Value = "'2.442"
$ws.Range("D26").Value = "'2.973"
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").Value = "'22.10"
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").Value = "'152.30"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").Value = "'5.407"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").Value = "'138.01"
$ws.Range("E30").Value = "  -3.62%  "
$ws.Range("D31").Value = "'7.931"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("D32").Value = "'2.508"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "1.830.02"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("D34").Value = "'1.026"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("D35").Value = "'0.08090"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("D36").Value = "'6.748"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "'0.02929"
$ws.Range("E37").Value = "  -2.12%  "

# Rows 38/39: FraxShare now ranks above Algorand; refresh name/link/price/volume
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'10.75"
$ws.Range("E38").Value = "  -4.16%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2677"
$ws.Range("E39").Value = "  -2.91%  "

$ws.Range("D40").Value = "'0.09121"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "'0.7604"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").Value = "'13.42"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("D43").Value = "'1.422"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").Value = "'16.26"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "'0.6967"
$ws.Range("E45").Value = "  -2.22%  "
$ws.Range("D46").Value = "'2.470"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").Value = "'4.097"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "'0.9997"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "'0.08344"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").Value = "'134.91"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Value = "'1.224"
$ws.Range("E51").Value = "  -3.44%  "

